$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-7, 10-19, 23-25: only the computed mat_range score (column H) was re-ranked.
$ws.Cells.Item(2,8).Value = 13.25581603006527
$ws.Cells.Item(3,8).Value = 13.03564410204013
$ws.Cells.Item(4,8).Value = 8.387878449008936
$ws.Cells.Item(5,8).Value = 8.356292063322577
$ws.Cells.Item(6,8).Value = 8.180015286402934
$ws.Cells.Item(7,8).Value = 5.393336665672788

# Row 8 and row 9 (ranks 7 & 8 for females) swapped identities after re-ranking.
$ws.Cells.Item(8,4).Value = 30
$ws.Cells.Item(8,5).Value = "60d5775a99b502eec8cf56b4"
$ws.Cells.Item(8,6).Value = "Shadaisia"
$ws.Cells.Item(8,8).Value = 5.339669197139461
$ws.Cells.Item(8,9).Value = "Black or African American"

$ws.Cells.Item(9,4).Value = 32
$ws.Cells.Item(9,5).Value = "6036f9b3b1842f8b659b18c7"
$ws.Cells.Item(9,6).Value = "Kellie"
$ws.Cells.Item(9,8).Value = 5.108019693417147
$ws.Cells.Item(9,9).Value = "White"

$ws.Cells.Item(10,8).Value = 4.268640122598316
$ws.Cells.Item(11,8).Value = 4.254495598246366
$ws.Cells.Item(12,8).Value = 2.420025270519735
$ws.Cells.Item(13,8).Value = 1.496024677253027
$ws.Cells.Item(14,8).Value = 14.35604799398173
$ws.Cells.Item(15,8).Value = 13.17756464437572
$ws.Cells.Item(16,8).Value = 8.158367614863963
$ws.Cells.Item(17,8).Value = 7.132419507397405
$ws.Cells.Item(18,8).Value = 6.306267974076017
$ws.Cells.Item(19,8).Value = 6.305467982787811

# Rows 20, 21, 22 (ranks 7, 8, 9 for males) rotated identities after re-ranking.
$ws.Cells.Item(20,4).Value = 30
$ws.Cells.Item(20,5).Value = "60c2341fe95d71ee52c043f0"
$ws.Cells.Item(20,6).Value = "Matthew"
$ws.Cells.Item(20,8).Value = 5.443833869706829
$ws.Cells.Item(20,9).Value = "White"

$ws.Cells.Item(21,4).Value = 32
$ws.Cells.Item(21,5).Value = "60bf9943e4e04642d4634ecc"
$ws.Cells.Item(21,6).Value = "Jamarii"
$ws.Cells.Item(21,8).Value = 5.210446373867417
$ws.Cells.Item(21,9).Value = "Black or African American"

$ws.Cells.Item(22,4).Value = 33
$ws.Cells.Item(22,5).Value = "60b322994d0b901954690036"
$ws.Cells.Item(22,6).Value = "Brennan"
$ws.Cells.Item(22,8).Value = 5.049269166493271

$ws.Cells.Item(23,8).Value = 3.295791998891051
$ws.Cells.Item(24,8).Value = 1.043195326962711
$ws.Cells.Item(25,8).Value = 0.2516358054655306
